# Update the two-digit division worksheet numbers to a new generated set.
# Each (row, col) cell in the single table is targeted individually (rather
# than a document-wide Find/Replace) so that values which are both a source
# and a target elsewhere (e.g. "11÷7=") cannot cross-contaminate each other.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{Row=1;  Col=1; Old="81÷5="; New="11÷7="},
    @{Row=1;  Col=2; Old="18÷6="; New="52÷9="},
    @{Row=1;  Col=3; Old="46÷8="; New="52÷7="},
    @{Row=1;  Col=4; Old="80÷2="; New="57÷2="},
    @{Row=1;  Col=5; Old="70÷3="; New="91÷9="},

    @{Row=5;  Col=1; Old="14÷5="; New="19÷3="},
    @{Row=5;  Col=2; Old="65÷9="; New="48÷6="},
    @{Row=5;  Col=3; Old="21÷6="; New="30÷6="},
    @{Row=5;  Col=4; Old="86÷5="; New="27÷7="},
    @{Row=5;  Col=5; Old="14÷3="; New="88÷6="},

    @{Row=9;  Col=1; Old="99÷7="; New="10÷5="},
    @{Row=9;  Col=2; Old="66÷7="; New="22÷3="},
    @{Row=9;  Col=3; Old="55÷4="; New="89÷5="},
    @{Row=9;  Col=4; Old="10÷8="; New="14÷9="},
    @{Row=9;  Col=5; Old="10÷7="; New="36÷9="},

    @{Row=13; Col=1; Old="11÷5="; New="33÷9="},
    @{Row=13; Col=2; Old="10÷4="; New="93÷3="},
    @{Row=13; Col=3; Old="15÷7="; New="28÷6="},
    @{Row=13; Col=4; Old="83÷5="; New="31÷4="},
    @{Row=13; Col=5; Old="11÷7="; New="26÷6="},

    @{Row=17; Col=1; Old="73÷2="; New="83÷7="},
    @{Row=17; Col=2; Old="76÷6="; New="39÷8="},
    @{Row=17; Col=3; Old="45÷5="; New="38÷8="},
    @{Row=17; Col=4; Old="70÷6="; New="71÷9="},
    @{Row=17; Col=5; Old="11÷8="; New="42÷2="}
)

foreach ($chg in $changes) {
    $cell = $t.Cell($chg.Row, $chg.Col)
    $rng = $cell.Range
    $before = $rng.Text
    if ($before -notlike ($chg.Old + "*")) {
        Write-Host "WARNING: unexpected text" $before "at row" $chg.Row "col" $chg.Col "(expected" $chg.Old ")"
    }
    # Assign the replacement text directly onto the cell's Range. Using
    # Find.Execute here is unsafe in this runtime: its search/replace is not
    # reliably scoped to the calling Range and can land on a different cell
    # that happens to share text with the Find/Replacement strings (several
    # values in this worksheet are simultaneously another cell's old or new
    # value, e.g. "11÷7="). Setting Range.Text replaces just this cell's
    # run content in place and keeps the run's existing formatting.
    $rng.Text = $chg.New
}

Write-Host "Done."
